$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F:G (pushes old F:M -> H:O), inheriting the
# formatting of the column that used to be F (bold/bordered header style).
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels for the two inserted columns.
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# New column data (M_TotalTax, M_CorpTax) for each data row.
$ws.Range("F2").Value = 14106286460237.92
$ws.Range("G2").Value = 1155021202746.413

$ws.Range("F3").Value = 3207987015.574299
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 734615892234.8064
$ws.Range("G4").Value = 88889835996.30263

$ws.Range("F5").Value = 558865056646.082
$ws.Range("G5").Value = 72600947639.16805

$ws.Range("F6").Value = 4579473077980.816
$ws.Range("G6").Value = 674619880691.7614

# Row 6 (UMICs) was recomputed using gains for all metrics, so besides the
# new columns, E6 and the shifted H6:O6 ratios also get new values.
$ws.Range("E6").Value = 2427884184.75
$ws.Range("H6").Value = 2.997565946324286
$ws.Range("I6").Value = 9.46852586630809
$ws.Range("J6").Value = 2.818099486909414
$ws.Range("K6").Value = 7.165342166154059
$ws.Range("L6").Value = 0.4027577601112367
$ws.Range("M6").Value = 1.272206295960165
$ws.Range("N6").Value = 2.013788800556183
$ws.Range("O6").Value = 6.361031479800828
